$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlink relationships pointing at the old column D
$ws.Range("D3").Hyperlinks.Delete()
$ws.Range("D5").Hyperlinks.Delete()

# Insert a new column before C: B(locator) C(action) D(value) -> B, C(new), D(action), E(value)
$ws.Columns("C:C").Insert()

# Approximate column C width to visually match column B (exact width unattainable via COM)
$ws.Columns("C:C").ColumnWidth = 15

# Header row
$ws.Range("B1").Value = "locatorType"
$ws.Range("C1").Value = "locatorValue"

# Row 2 (open browser) - locator was NA -> split into NA / NA
$ws.Range("B2").Value = "NA"
$ws.Range("C2").Value = "NA"

# Row 3 (launch url) - locator was NA -> split into NA / NA
$ws.Range("B3").Value = "NA"
$ws.Range("C3").Value = "NA"

# Row 4 (enter email address) - locator was "name = username" -> "name  " / "username"
$ws.Range("B4").Value = "name  "
$ws.Range("C4").Value = "username"

# Row 5 (enter password) - locator was "name = password" -> "name  " / "password"
$ws.Range("B5").Value = "name  "
$ws.Range("C5").Value = "password"

# Re-add hyperlinks at their new location (column E now) and restore original hyperlink style
$ws.Hyperlinks.Add($ws.Range("E3"), "https://classic.crmpro.com/index.html")
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:Chaithu@17")
$ws.Range("E3").Style = "Hyperlink"
$ws.Range("E5").Style = "Hyperlink"

# Update selection to C1
$ws.Range("C1").Select()
